# Update countries & provincias Spain
# Applies the data refresh described in the commit: updated case counts for
# several countries, a handful of countries swapping rank (because their
# updated totals changed their sort order), and a refreshed "last updated"
# timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 14:41"

# --- Helper to rewrite a whole data row (country + 7 numeric columns) ----
function Set-Row([int]$Row, [string]$Country, [double]$CasosTotales, [double]$NuevosCasos, [double]$CasosActivos, [double]$Recuperados, [double]$CasosCriticos, [double]$MuertesHoy, [double]$Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Estados Unidos
Set-Row 4 "Estados Unidos" 2837612 423 1191838 1514271 0 18 131503

# Brasil
Set-Row 5 "Brasil" 1502424 1071 916147 524232 0 55 62045

# Alemania
Set-Row 18 "Alemania" 196738 21 181000 6674 0 0 9064

# Dinamarca
Set-Row 67 "Dinamarca" 12832 17 11817 409 0 0 606

# Uzbekistan
Set-Row 72 "Uzbekistan" 9326 248 6161 3137 0 1 28

# Kenia & Senegal swap rank (row 79 / row 80)
Set-Row 79 "Kenia" 7188 247 2109 4927 0 0 152
Set-Row 80 "Senegal" 7164 110 4666 2373 0 4 125

# Republica de Yibuti
Set-Row 91 "Republica de Yibuti" 4736 21 4580 101 0 0 55

# Croacia & Somalia swap rank (row 102 / row 103)
Set-Row 102 "Croacia" 3008 96 2168 728 0 2 112
Set-Row 103 "Somalia" 2944 0 951 1903 0 0 90

# Surinam
Set-Row 154 "Surinam" 547 0 254 280 0 0 13

# Dominica & Fiyi swap rank (row 205 / row 206) - data tied, only name order changes
Set-Row 205 "Dominica" 18 0 18 0 0 0 0
Set-Row 206 "Fiyi" 18 0 18 0 0 0 0

# Islas Malvinas & Groenlandia swap rank (row 209 / row 210) - data tied, only name order changes
Set-Row 209 "Islas Malvinas" 13 0 13 0 0 0 0
Set-Row 210 "Groenlandia" 13 0 13 0 0 0 0
